# "Various changes to update to VS26"
# Adds a new row of device data (Waldorf Iridium Desktop) to the table on
# Sheet1, and updates the active selection/scroll position to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Populate the previously-empty row 39 with the new device entry.
$ws.Range("A39").Value = "Waldorf Iridium Desktop"
$ws.Range("B39").Value = "Synthesizer"
$ws.Range("C39").Value = "KS"
$ws.Range("D39").Value = "Iridium"
$ws.Range("E39").Value = "(bidirectional FB) Synth"
$ws.Range("F39").Value = "(bidirectional FB) Synth"
$ws.Range("G39").Value = "Iridium (MIDI 2.0)"
$ws.Range("H39").Value = "Iridium (MIDI 2.0)"
$ws.Range("I39").Value = "Using MIDI 2.0 firmware."

# Scroll the sheet so row 9 is at the top, then select A40 as the active
# cell, matching the saved view state.
$excel.Goto($ws.Range("A9"), $true)
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A40").Select()
